# New simulation files for schemes report
# Updates the "UniformF-HW05" sheet:
#  - relabels the existing HKL-scheme rows (B3:B19) with the new scheme names
#  - relabels the [h,k,l] column headers (C2:T2) to their new order
#  - appends 10 new HKL-scheme rows (20-29) with the same layout as the
#    existing data rows (A: index, B: scheme name, C:T: all 1's)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Relabel existing scheme-name column (B3:B19) with the new names
# ---------------------------------------------------------------------
$newSchemeNames = @(
    "Spiral5",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Thomas Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD"
)

for ($i = 0; $i -lt $newSchemeNames.Count; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 2).Value = $newSchemeNames[$i]
}

# ---------------------------------------------------------------------
# 2. Relabel the [h,k,l] header row (C2:T2) with the new ordering
# ---------------------------------------------------------------------
$newHeaders = @(
    "[3, 2, 1]",
    "[1, 1, 0]",
    "[3, 1, 0]",
    "[2, 2, 2]",
    "[2, 0, 0]",
    "[2, 2, 0]",
    "[2, 1, 1]",
    "[4, 0, 0]",
    "1Pair-A",
    "1Pair-B",
    "2Pairs-A",
    "2Pairs-B",
    "3Pairs-A",
    "3Pairs-B",
    "3Pairs-C",
    "4Pairs",
    "5A4F",
    "MaxUnique"
)

for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $col = 3 + $i
    $ws.Cells.Item(2, $col).Value = $newHeaders[$i]
}

# ---------------------------------------------------------------------
# 3. Append 10 new data rows (20-29), matching the existing table layout
# ---------------------------------------------------------------------
$newRowNames = @(
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex"
)

# Copy the formatting (styles) of an existing full data row (row 3) down
# onto each new row before writing values, so the new rows pick up the
# same bold/bordered index-column style used by the rest of the table.
$formatSource = $ws.Range("A3:T3")

for ($i = 0; $i -lt $newRowNames.Count; $i++) {
    $row = 20 + $i
    $index = 18 + $i

    $formatSource.Copy()
    $ws.Range("A" + $row + ":T" + $row).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $index
    $ws.Cells.Item($row, 2).Value = $newRowNames[$i]

    for ($col = 3; $col -le 20; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}

$excel.CutCopyMode = 0
